# Applies the "output generated at 456a3b4" update to 北京-漫展信息.xlsx
# Sheets: 1=展览 (Exhibition), 2=演出 (Performance), 3=本地生活 (Local life),
#         4=全部类型 (All types, a combined/sorted view of the first three)

$wb = $excel.ActiveWorkbook

$wsExpo   = $wb.Worksheets.Item("展览")
$wsShow   = $wb.Worksheets.Item("演出")
$wsLocal  = $wb.Worksheets.Item("本地生活")
$wsAll    = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------------------
# Sheet 1: 展览 -- "想去人数" (F) / "最低票价" (G) counter bumps
# ---------------------------------------------------------------------------
$wsExpo.Range("F3").Value  = 429
$wsExpo.Range("F5").Value  = 3948
$wsExpo.Range("G7").Value  = 80
$wsExpo.Range("F9").Value  = 3221
$wsExpo.Range("G9").Value  = 80
$wsExpo.Range("F10").Value = 536
$wsExpo.Range("G10").Value = 75
$wsExpo.Range("G11").Value = 65
$wsExpo.Range("F15").Value = 472
$wsExpo.Range("F16").Value = 23
$wsExpo.Range("F20").Value = 316
$wsExpo.Range("F21").Value = 465
$wsExpo.Range("F24").Value = 345
$wsExpo.Range("F25").Value = 14
$wsExpo.Range("F28").Value = 174
$wsExpo.Range("F31").Value = 4459
$wsExpo.Range("F32").Value = 4362
$wsExpo.Range("F34").Value = 307
$wsExpo.Range("F35").Value = 77
$wsExpo.Range("F36").Value = 18
$wsExpo.Range("F37").Value = 1171
$wsExpo.Range("F40").Value = 506
$wsExpo.Range("F42").Value = 1326
$wsExpo.Range("F43").Value = 188
$wsExpo.Range("F44").Value = 139
$wsExpo.Range("F48").Value = 68

# ---------------------------------------------------------------------------
# Sheet 2: 演出 -- row 4 ("最后的莫西干人...") is no longer on sale
# ---------------------------------------------------------------------------
$wsShow.Range("G4").Value = "不可售"

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 -- "想去人数" (F) counter bumps
# ---------------------------------------------------------------------------
$wsLocal.Range("F4").Value = 2344
$wsLocal.Range("F5").Value = 57

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 -- combined/sorted feed. A new event
# (" 北京·万游引力嘉年华 配音演员赵成晨&尘霜满眸 广播剧《奕曲同工》专场见面&签售会")
# now sorts into the same 2024-08-24 slot as rows 7-10, so those rows'
# contents shift down by one, and the previously-last same-day row
# ("北京·最后的莫西干人...") drops off this combined view (it now shows as
# "不可售" over on the 演出 sheet instead). Net row count is unchanged, so
# this is expressed as in-place cell overwrites rather than a real insert.
# ---------------------------------------------------------------------------
$wsAll.Range("F4").Value = 429
$wsAll.Range("F6").Value = 3948

$wsAll.Range("C7").Value = " 北京·万游引力嘉年华 配音演员赵成晨&尘霜满眸 广播剧《奕曲同工》专场见面&签售会"
$wsAll.Range("D7").Value = "金蝉西路甲1号（地铁七号线南楼梓庄站） 北京酷车国际汇展中心"
$wsAll.Range("E7").Value = "2024.08.24 11:00-08.24 17:00"
$wsAll.Range("F7").Value = 225
$wsAll.Range("G7").Value = 288
$wsAll.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=89054"
$wsAll.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202407/FadWpN3x1720599868028.jpeg"

$wsAll.Range("C8").Value = "北京·ACY动漫游戏展1st"
$wsAll.Range("D8").Value = "崇文门外大街18号 北京国瑞购物中心"
$wsAll.Range("E8").Value = "2024.08.24 10:00-08.25 17:00"
$wsAll.Range("F8").Value = 2602
$wsAll.Range("G8").Value = 80
$wsAll.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=87851"
$wsAll.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202406/jKUUlXAR1718967902236.jpeg"

$wsAll.Range("C9").Value = "北京·“不健全关系”专题聚会【免票活动】"
$wsAll.Range("D9").Value = "王府井大街88号 北京王府井银泰in88购物中心"
$wsAll.Range("E9").Value = "2024.08.24 14:00-08.24 18:00"
$wsAll.Range("F9").Value = 85
$wsAll.Range("G9").Value = 50
$wsAll.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=90562"
$wsAll.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202408/qBI8D5Ji1723624479890.jpeg"

$wsAll.Range("C10").Value = "北京·万游引力夏日动漫游戏狂欢节"
$wsAll.Range("D10").Value = "金蝉西路甲1号（地铁七号线南楼梓庄站） 北京酷车国际汇展中心"
$wsAll.Range("E10").Value = "2024.08.24 10:00-08.25 17:00"
$wsAll.Range("F10").Value = 3221
$wsAll.Range("G10").Value = 80
$wsAll.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=83880"
$wsAll.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202407/3EF1Am6T1720430616435.jpeg"

$wsAll.Range("F11").Value = 536
$wsAll.Range("G11").Value = 75

$wsAll.Range("G12").Value = 65

$wsAll.Range("F16").Value = 472
$wsAll.Range("F17").Value = 23
$wsAll.Range("F29").Value = 4459
$wsAll.Range("F30").Value = 4362
$wsAll.Range("F32").Value = 18
$wsAll.Range("F33").Value = 1171
$wsAll.Range("F38").Value = 506
$wsAll.Range("F43").Value = 1326
$wsAll.Range("F44").Value = 188
$wsAll.Range("F48").Value = 68
